$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, date range) ---
$ws.Range("A8").Value = "Volume 31   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/24/2024  Through  6/30/2024"

# --- Data cell updates (rows 14-31) ---
$ws.Range("N14").Value = -72.727272727272
$ws.Range("F15").Value = 2
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("N15").Value = -54.545454545454
$ws.Range("D16").Value = 3
$ws.Range("F16").Value = 15
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 99
$ws.Range("J16").Value = 89
$ws.Range("L16").Value = 17.857142857142
$ws.Range("M16").Value = -40.718562874251
$ws.Range("N16").Value = -83.928571428571
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -69.230769230769
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = 6.896551724137
$ws.Range("I17").Value = 154
$ws.Range("J17").Value = 154
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 16.666666666666
$ws.Range("M17").Value = 67.391304347826
$ws.Range("N17").Value = -40.077821011673
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 29
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = 26.086956521739
$ws.Range("I18").Value = 176
$ws.Range("J18").Value = 111
$ws.Range("K18").Value = 58.558558558558
$ws.Range("L18").Value = 15.78947368421
$ws.Range("M18").Value = -7.853403141361
$ws.Range("N18").Value = -72.01907790143
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = -38.095238095238
$ws.Range("F19").Value = 48
$ws.Range("G19").Value = 80
$ws.Range("H19").Value = -40
$ws.Range("I19").Value = 353
$ws.Range("J19").Value = 374
$ws.Range("K19").Value = -5.614973262032
$ws.Range("L19").Value = 15.737704918032
$ws.Range("M19").Value = 64.186046511627
$ws.Range("N19").Value = 38.976377952755
$ws.Range("C20").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D20").Value = 6
$ws.Range("C16").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = -83.333333333333
$ws.Range("K16").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = -46.153846153846
$ws.Range("I20").Value = 63
$ws.Range("J20").Value = 77
$ws.Range("K20").Value = -18.181818181818
$ws.Range("L20").Value = -10
$ws.Range("M20").Value = -16
$ws.Range("N20").Value = -85.245901639344
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 53
$ws.Range("E21").Value = -50.943396226415
$ws.Range("F21").Value = 132
$ws.Range("G21").Value = 155
$ws.Range("H21").Value = -14.838709677419
$ws.Range("I21").Value = 858
$ws.Range("J21").Value = 816
$ws.Range("K21").Value = 5.147058823529
$ws.Range("L21").Value = 14.552736982643
$ws.Range("M21").Value = 15.322580645161
$ws.Range("N21").Value = -61.281588447653
$ws.Range("D22").Value = 1
$ws.Range("C16").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("K16").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("G22").Value = 2
$ws.Range("J22").Value = 16
$ws.Range("K22").Value = -25
$ws.Range("L22").Value = -50
$ws.Range("M22").Value = -45.454545454545
$ws.Range("C23").Value = 1
$ws.Range("E23").Value = -83.333333333333
$ws.Range("F23").Value = 11
$ws.Range("G23").Value = 22
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 89
$ws.Range("J23").Value = 101
$ws.Range("K23").Value = -11.881188118811
$ws.Range("L23").Value = 17.105263157894
$ws.Range("M23").Value = 34.848484848484
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 18.75
$ws.Range("F24").Value = 107
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = 28.915662650602
$ws.Range("I24").Value = 547
$ws.Range("J24").Value = 519
$ws.Range("K24").Value = 5.394990366088
$ws.Range("L24").Value = -10.474631751227
$ws.Range("M24").Value = -6.972789115646
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 75
$ws.Range("F25").Value = 23
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = 64.285714285714
$ws.Range("I25").Value = 185
$ws.Range("J25").Value = 67
$ws.Range("K25").Value = 176.119402985075
$ws.Range("L25").Value = 28.472222222222
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -72.727272727272
$ws.Range("F26").Value = 40
$ws.Range("G26").Value = 48
$ws.Range("H26").Value = -16.666666666666
$ws.Range("I26").Value = 278
$ws.Range("J26").Value = 234
$ws.Range("K26").Value = 18.803418803418
$ws.Range("L26").Value = 6.106870229007
$ws.Range("M26").Value = 14.40329218107
$ws.Range("F27").Value = 3
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -66.666666666666
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 20
$ws.Range("I28").Value = 29
$ws.Range("J28").Value = 34
$ws.Range("K28").Value = -14.705882352941
$ws.Range("L28").Value = 20.833333333333
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("H29").Value = -100
$ws.Range("N29").Value = -81.818181818181
$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("H30").Value = -100
$ws.Range("N30").Value = -84.375
$ws.Range("F31").Value = 1
$ws.Range("C16").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = -66.666666666666
$ws.Range("I31").Value = 16
$ws.Range("J31").Value = 7
$ws.Range("K31").Value = 128.571428571429
$ws.Range("L31").Value = 33.333333333333

# --- Donor cells whose own values also change; updated last ---
$ws.Range("C16").Value = 3
$ws.Range("K16").Value = 11.235955056179

$ws.Range("A1").Select()
$excel.CutCopyMode = $false
